# Refresh the crypto price-tracker sheet with the latest scrape results
# (GitHub Actions cron job). Prices/volume figures are stored as plain text
# in this sheet (e.g. "37.450.10", "  +2.25%  "), and a handful of rows were
# also re-ranked, so a few Coin/Link cells swap with their neighbours too.
#
# AsText marks price cells whose new value happens to look like a plain
# number (e.g. "0.995", "251.90"). Excel's COM `.Value` setter auto-converts
# such strings to a numeric cell, so for those we prefix a literal quote
# (exactly like typing `'0.995` into Excel) to force text, then reset the
# resulting quote-prefix cell style back to Normal so no stray formatting
# sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '37.450.10'; AsText = $false },
    @{ Cell = 'E2'; Value = '  +2.25%  '; AsText = $false },
    @{ Cell = 'D3'; Value = '2.083.88'; AsText = $false },
    @{ Cell = 'E3'; Value = '  -0.89%  '; AsText = $false },
    @{ Cell = 'D4'; Value = '0.995'; AsText = $true },
    @{ Cell = 'E4'; Value = '  -0.45%  '; AsText = $false },
    @{ Cell = 'D5'; Value = '251.90'; AsText = $true },
    @{ Cell = 'E5'; Value = '  +1.44%  '; AsText = $false },
    @{ Cell = 'D6'; Value = '0.665'; AsText = $true },
    @{ Cell = 'E6'; Value = '  -0.30%  '; AsText = $false },
    @{ Cell = 'D7'; Value = '0.999'; AsText = $true },
    @{ Cell = 'E7'; Value = '  -0.05%  '; AsText = $false },
    @{ Cell = 'D8'; Value = '56.23'; AsText = $true },
    @{ Cell = 'E8'; Value = '  +23.63%  '; AsText = $false },
    @{ Cell = 'D9'; Value = '62.28'; AsText = $true },
    @{ Cell = 'E9'; Value = '  +2.73%  '; AsText = $false },
    @{ Cell = 'D10'; Value = '0.381'; AsText = $true },
    @{ Cell = 'E10'; Value = '  +4.17%  '; AsText = $false },
    @{ Cell = 'D11'; Value = '0.0754'; AsText = $true },
    @{ Cell = 'E11'; Value = '  +3.57%  '; AsText = $false },
    @{ Cell = 'E12'; Value = '  +7.37%  '; AsText = $false },
    @{ Cell = 'D13'; Value = '15.41'; AsText = $true },
    @{ Cell = 'E13'; Value = '  +5.35%  '; AsText = $false },
    @{ Cell = 'D14'; Value = '2.388.66'; AsText = $false },
    @{ Cell = 'E14'; Value = '  -0.26%  '; AsText = $false },
    @{ Cell = 'D15'; Value = '0.856'; AsText = $true },
    @{ Cell = 'E15'; Value = '  +2.53%  '; AsText = $false },
    @{ Cell = 'E16'; Value = '  +6.07%  '; AsText = $false },
    @{ Cell = 'D17'; Value = '2.084.85'; AsText = $false },
    @{ Cell = 'E17'; Value = '  -0.35%  '; AsText = $false },
    @{ Cell = 'D18'; Value = '37.206.44'; AsText = $false },
    @{ Cell = 'E18'; Value = '  +1.56%  '; AsText = $false },
    @{ Cell = 'D19'; Value = '73.33'; AsText = $true },
    @{ Cell = 'E19'; Value = '  +1.74%  '; AsText = $false },
    @{ Cell = 'D20'; Value = '14.59'; AsText = $true },
    @{ Cell = 'E20'; Value = '  +13.87%  '; AsText = $false },
    @{ Cell = 'D21'; Value = '0.0₃0854'; AsText = $false },
    @{ Cell = 'D22'; Value = '241.26'; AsText = $true },
    @{ Cell = 'E22'; Value = '  +0.57%  '; AsText = $false },
    @{ Cell = 'D23'; Value = '5.28'; AsText = $true },
    @{ Cell = 'E23'; Value = '  +5.95%  '; AsText = $false },
    @{ Cell = 'E24'; Value = '  -0.07%  '; AsText = $false },
    @{ Cell = 'E25'; Value = '  +0.12%  '; AsText = $false },
    @{ Cell = 'D26'; Value = '171.79'; AsText = $true },
    @{ Cell = 'E26'; Value = '  +1.17%  '; AsText = $false },
    @{ Cell = 'D27'; Value = '9.28'; AsText = $true },
    @{ Cell = 'E27'; Value = '  +4.90%  '; AsText = $false },
    @{ Cell = 'D28'; Value = '20.89'; AsText = $true },
    @{ Cell = 'E28'; Value = '  +0.70%  '; AsText = $false },
    @{ Cell = 'E29'; Value = '  +3.01%  '; AsText = $false },
    @{ Cell = 'E30'; Value = '  +1.56%  '; AsText = $false },
    @{ Cell = 'D31'; Value = '23.58'; AsText = $true },
    @{ Cell = 'E31'; Value = '  +6.59%  '; AsText = $false },
    @{ Cell = 'D32'; Value = '1.10'; AsText = $true },
    @{ Cell = 'E32'; Value = '  +21.27%  '; AsText = $false },
    @{ Cell = 'D33'; Value = '4.56'; AsText = $true },
    @{ Cell = 'E33'; Value = '  +3.76%  '; AsText = $false },
    @{ Cell = 'E34'; Value = '  +6.44%  '; AsText = $false },
    @{ Cell = 'D35'; Value = '0.0905'; AsText = $true },
    @{ Cell = 'E35'; Value = '  +0.56%  '; AsText = $false },
    @{ Cell = 'D36'; Value = '4.31'; AsText = $true },
    @{ Cell = 'E36'; Value = '  +7.16%  '; AsText = $false },
    @{ Cell = 'D37'; Value = '0.996'; AsText = $true },
    @{ Cell = 'E37'; Value = '  -0.39%  '; AsText = $false },
    @{ Cell = 'B38'; Value = 'LidoDAOToken'; AsText = $false },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; AsText = $false },
    @{ Cell = 'D38'; Value = '2.27'; AsText = $true },
    @{ Cell = 'E38'; Value = '  -2.25%  '; AsText = $false },
    @{ Cell = 'B39'; Value = 'WEMIXToken'; AsText = $false },
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; AsText = $false },
    @{ Cell = 'D39'; Value = '1.84'; AsText = $true },
    @{ Cell = 'E39'; Value = '  -2.95%  '; AsText = $false },
    @{ Cell = 'E40'; Value = '  +0.11%  '; AsText = $false },
    @{ Cell = 'E41'; Value = '  +5.12%  '; AsText = $false },
    @{ Cell = 'D42'; Value = '17.91'; AsText = $true },
    @{ Cell = 'E42'; Value = '  +12.12%  '; AsText = $false },
    @{ Cell = 'D43'; Value = '1.18'; AsText = $true },
    @{ Cell = 'E43'; Value = '  +1.39%  '; AsText = $false },
    @{ Cell = 'D44'; Value = '0.0981'; AsText = $true },
    @{ Cell = 'E44'; Value = '  +18.73%  '; AsText = $false },
    @{ Cell = 'D45'; Value = '99.99'; AsText = $true },
    @{ Cell = 'E45'; Value = '  +1.55%  '; AsText = $false },
    @{ Cell = 'B46'; Value = 'FTXToken'; AsText = $false },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'; AsText = $false },
    @{ Cell = 'D46'; Value = '4.32'; AsText = $true },
    @{ Cell = 'E46'; Value = '  +118.15%  '; AsText = $false },
    @{ Cell = 'B47'; Value = 'HuobiToken'; AsText = $false },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; AsText = $false },
    @{ Cell = 'D47'; Value = '2.81'; AsText = $true },
    @{ Cell = 'E47'; Value = '  +0.43%  '; AsText = $false },
    @{ Cell = 'D48'; Value = '1.330.52'; AsText = $false },
    @{ Cell = 'E48'; Value = '  -1.26%  '; AsText = $false },
    @{ Cell = 'E49'; Value = '  +3.85%  '; AsText = $false },
    @{ Cell = 'B50'; Value = 'FraxShare'; AsText = $false },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; AsText = $false },
    @{ Cell = 'D50'; Value = '7.03'; AsText = $true },
    @{ Cell = 'E50'; Value = '  +11.84%  '; AsText = $false },
    @{ Cell = 'B51'; Value = 'RenderToken'; AsText = $false },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; AsText = $false },
    @{ Cell = 'D51'; Value = '2.35'; AsText = $true },
    @{ Cell = 'E51'; Value = '  +5.30%  '; AsText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.AsText) {
        $range.Value = "'" + $u.Value
        $range.Style = 'Normal'
    } else {
        $range.Value = $u.Value
    }
}
